$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1398.375
$ws.Range("I9").Value = 997.8333
$ws.Range("K9").Value = 997.8333
$ws.Range("M9").Value = -828.8333
$ws.Range("H29").Value = 5062
$ws.Range("J29").Value = 5062
$ws.Range("L29").Value = 15186
$ws.Range("N29").Value = -15748
$ws.Range("H33").Value = 695.8946999999999
$ws.Range("I33").Value = 110.411766
$ws.Range("J33").Value = 5672.5
$ws.Range("K33").Value = 110.411766
$ws.Range("L33").Value = 5672.5
$ws.Range("M33").Value = 118.588234
$ws.Range("N33").Value = -6130.5
$ws.Range("H38").Value = 10401.4
$ws.Range("I38").Value = 1004
$ws.Range("J38").Value = 16666.334
$ws.Range("K38").Value = 3012
$ws.Range("L38").Value = 49999.00199999999
$ws.Range("M38").Value = -2640
$ws.Range("N38").Value = -50743.00199999999
$ws.Range("H43").Value = 5125
$ws.Range("J43").Value = 6166.6665
$ws.Range("L43").Value = 6166.6665
$ws.Range("N43").Value = -6304.6665
$ws.Range("H58").Value = 3017.25
$ws.Range("J58").Value = 5999.5
$ws.Range("L58").Value = 17998.5
$ws.Range("N58").Value = -18298.5
$ws.Range("H69").Value = 4665.6665
$ws.Range("J69").Value = 5000
$ws.Range("L69").Value = 15000
$ws.Range("N69").Value = -16748
$ws.Range("H72").Value = 4665.6665
$ws.Range("J72").Value = 5000
$ws.Range("L72").Value = 45000
$ws.Range("N72").Value = -53736
$ws.Range("H112").Value = 2955.077
$ws.Range("J112").Value = 3083
$ws.Range("L112").Value = 9249
$ws.Range("N112").Value = -11465
$ws.Range("H121").Value = 481.83334
$ws.Range("J121").Value = 481.83334
$ws.Range("L121").Value = 1445.50002
$ws.Range("N121").Value = -4939.500019999999
$ws.Range("H132").Value = 1815.8422
$ws.Range("I132").Value = 1749.0646
$ws.Range("J132").Value = 2111.5715
$ws.Range("K132").Value = 5247.1938
$ws.Range("L132").Value = 6334.7145
$ws.Range("M132").Value = -2717.1938
$ws.Range("N132").Value = -11394.7145
$ws.Range("H137").Value = 3296.7742
$ws.Range("I137").Value = 1488.7778
$ws.Range("K137").Value = 4466.3334
$ws.Range("M137").Value = -1916.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 4857
$ws.Range("I4").Value = 2042.5
$ws.Range("J4").Value = 6733.3335
$ws.Range("K4").Value = 2042.5
$ws.Range("L4").Value = 6733.3335
$ws.Range("M4").Value = -1926.5
$ws.Range("N4").Value = -6965.3335
$ws.Range("H32").Value = 15366.493
$ws.Range("I32").Value = 6245.65
$ws.Range("K32").Value = 6245.65
$ws.Range("M32").Value = -5958.65
$ws.Range("H74").Value = 3522.0667
$ws.Range("I74").Value = 1020.375
$ws.Range("J74").Value = 6381.143
$ws.Range("K74").Value = 1020.375
$ws.Range("L74").Value = 6381.143
$ws.Range("M74").Value = -146.375
$ws.Range("N74").Value = -8129.143
$ws.Range("H77").Value = 3522.0667
$ws.Range("I77").Value = 1020.375
$ws.Range("J77").Value = 6381.143
$ws.Range("K77").Value = 5101.875
$ws.Range("L77").Value = 31905.715
$ws.Range("M77").Value = -733.875
$ws.Range("N77").Value = -40641.715
$ws.Range("H122").Value = 2857.25
$ws.Range("I122").Value = 1976.375
$ws.Range("K122").Value = 5929.125
$ws.Range("M122").Value = -3479.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5520.1177
$ws.Range("I105").Value = 4195.25
$ws.Range("J105").Value = 8699.799999999999
$ws.Range("K105").Value = 4195.25
$ws.Range("L105").Value = 8699.799999999999
$ws.Range("M105").Value = -2448.25
$ws.Range("N105").Value = -12193.8
$ws.Range("H134").Value = 2296.4583
$ws.Range("I134").Value = 1576.5714
$ws.Range("J134").Value = 7335.6665
$ws.Range("K134").Value = 4729.7142
$ws.Range("L134").Value = 22006.9995
$ws.Range("M134").Value = -2194.7142
$ws.Range("N134").Value = -27076.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4152.077
$ws.Range("I31").Value = 2822.3076
$ws.Range("J31").Value = 5481.846
$ws.Range("K31").Value = 2822.3076
$ws.Range("L31").Value = 5481.846
$ws.Range("M31").Value = -2527.3076
$ws.Range("N31").Value = -6071.846
$ws.Range("H34").Value = 4152.077
$ws.Range("I34").Value = 2822.3076
$ws.Range("J34").Value = 5481.846
$ws.Range("K34").Value = 2822.3076
$ws.Range("L34").Value = 5481.846
$ws.Range("M34").Value = -2620.3076
$ws.Range("N34").Value = -5885.846
$ws.Range("H58").Value = 9072.333000000001
$ws.Range("I58").Value = 8887.5
$ws.Range("J58").Value = 9164.75
$ws.Range("K58").Value = 8887.5
$ws.Range("L58").Value = 9164.75
$ws.Range("M58").Value = -8684.5
$ws.Range("N58").Value = -9570.75
$ws.Range("H105").Value = 2595.2
$ws.Range("I105").Value = 2577.3333
$ws.Range("J105").Value = 2622
$ws.Range("K105").Value = 2577.3333
$ws.Range("L105").Value = 2622
$ws.Range("M105").Value = -830.3332999999998
$ws.Range("N105").Value = -6116
$ws.Range("H136").Value = 9072.333000000001
$ws.Range("I136").Value = 8887.5
$ws.Range("J136").Value = 9164.75
$ws.Range("K136").Value = 26662.5
$ws.Range("L136").Value = 27494.25
$ws.Range("M136").Value = -24112.5
$ws.Range("N136").Value = -32594.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2501632
$ws.Range("I4").Value = 3333676.2
$ws.Range("J4").Value = 5499
$ws.Range("K4").Value = 10001028.6
$ws.Range("L4").Value = 16497
$ws.Range("M4").Value = -10000916.6
$ws.Range("N4").Value = -16721
$ws.Range("H80").Value = 5954.4443
$ws.Range("I80").Value = 5897.5
$ws.Range("K80").Value = 17692.5
$ws.Range("M80").Value = -16756.5
$ws.Range("H83").Value = 5954.4443
$ws.Range("I83").Value = 5897.5
$ws.Range("K83").Value = 53077.5
$ws.Range("M83").Value = -48397.5
$ws.Range("H92").Value = 425
$ws.Range("I92").Value = 425
$ws.Range("K92").Value = 1275
$ws.Range("M92").Value = -27
$ws.Range("H94").Value = 2111
$ws.Range("I94").Value = 2111
$ws.Range("K94").Value = 6333
$ws.Range("M94").Value = -5657

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6995
$ws.Range("J70").Value = 6995
$ws.Range("L70").Value = 6995
$ws.Range("N70").Value = -7535
$ws.Range("H73").Value = 6995
$ws.Range("J73").Value = 6995
$ws.Range("L73").Value = 6995
$ws.Range("N73").Value = -8867

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3123.3333
$ws.Range("I46").Value = 2497.7778
$ws.Range("K46").Value = 2497.7778
$ws.Range("M46").Value = -2309.7778
$ws.Range("H55").Value = 688.5714
$ws.Range("I55").Value = 355.5625
$ws.Range("J55").Value = 1754.2
$ws.Range("K55").Value = 355.5625
$ws.Range("L55").Value = 1754.2
$ws.Range("M55").Value = -182.5625
$ws.Range("N55").Value = -2100.2
$ws.Range("H61").Value = 5156.647
$ws.Range("J61").Value = 4210.75
$ws.Range("L61").Value = 4210.75
$ws.Range("N61").Value = -4614.75
$ws.Range("H82").Value = 2974.7778
$ws.Range("I82").Value = 2681.8572
$ws.Range("K82").Value = 2681.8572
$ws.Range("M82").Value = -2320.8572
$ws.Range("H85").Value = 2974.7778
$ws.Range("I85").Value = 2681.8572
$ws.Range("K85").Value = 2681.8572
$ws.Range("M85").Value = -1433.8572
$ws.Range("H93").Value = 2037
$ws.Range("I93").Value = 1555.5
$ws.Range("K93").Value = 1555.5
$ws.Range("M93").Value = -307.5
$ws.Range("H104").Value = 23081.25
$ws.Range("J104").Value = 23081.25
$ws.Range("L104").Value = 23081.25
$ws.Range("N104").Value = -30069.25
$ws.Range("H113").Value = 5156.647
$ws.Range("J113").Value = 4210.75
$ws.Range("L113").Value = 4210.75
$ws.Range("N113").Value = -8550.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H60").Value = 130000
$ws.Range("J60").Value = 150000
$ws.Range("L60").Value = 150000
$ws.Range("N60").Value = -151644
$ws.Range("H81").Value = 2853.6924
$ws.Range("I81").Value = 2854.3635
$ws.Range("K81").Value = 5708.727
$ws.Range("M81").Value = -4647.727
$ws.Range("H84").Value = 2853.6924
$ws.Range("I84").Value = 2854.3635
$ws.Range("K84").Value = 28543.635
$ws.Range("M84").Value = -23239.635
$ws.Range("H123").Value = 5000
$ws.Range("I123").Value = 5000
$ws.Range("K123").Value = 5000
$ws.Range("M123").Value = -100
